# Applies the TaskList.xlsx edit described by the commit:
# "Refactor code and fix formatting issues"
#
# Summary of content changes on sheet task_list_1:
#  - Clear stray "N" values that were mistakenly left in column D for rows
#    that have no actual Result yet (rows 3-13, 17-18, 37-48).
#  - Row 15: mark Result D15 = "Y" and Remarks E15 = "Failed".
#  - Row 16: fill in the Remarks (E16 = "Success") and the numeric results
#    that were missing (F,G,H,J,K,L,M,O), matching the pattern already
#    present on row 14.
#  - Move the active cell selection to F12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task_list_1")

# Rows whose column D held a leftover "N" that should be cleared.
$rowsToClear = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 17, 18, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48)
foreach ($r in $rowsToClear) {
    $ws.Range("D$r").Value = ""
}

# Row 15: Result -> Y, Remarks -> Failed
$ws.Range("D15").Value = "Y"
$ws.Range("E15").Value = "Failed"

# Row 16: fill in Remarks + numeric results (mirrors row 14's values)
$ws.Range("E16").Value = "Success"
$ws.Range("F16").Value = 2662
$ws.Range("G16").Value = 1767
$ws.Range("H16").Value = 0.3
$ws.Range("J16").Value = -69
$ws.Range("K16").Value = 4388
$ws.Range("L16").Value = 2914
$ws.Range("M16").Value = 0.49
$ws.Range("O16").Value = 86

# Update the active selection to F12
$ws.Range("F12").Select()
